$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 11)
$ws.Range("H11").Value = "Fault coverage "
$ws.Range("I11").Value = "Test coverage"
$ws.Range("J11").Value = "Clock cycles"

# Data row (row 12)
$ws.Range("G12").Value = "G08_MUL"
$ws.Range("H12").Value = 0.3028
$ws.Range("I12").Value = 0.3032

$headerRange = $ws.Range("G11:J12")
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$ws.Range("G11:I11").HorizontalAlignment = -4108
$ws.Range("G11:I11").VerticalAlignment = -4108
$ws.Range("G12").HorizontalAlignment = -4108
$ws.Range("G12").VerticalAlignment = -4108

$ws.Range("H12:I12").NumberFormat = "0.00%"
$ws.Range("H12:I12").HorizontalAlignment = -4108
$ws.Range("H12:I12").VerticalAlignment = -4108

$ws.Columns("H").ColumnWidth = 17.140625
$ws.Columns("I").ColumnWidth = 17
$ws.Columns("J").ColumnWidth = 14.140625

$ws.Range("J17").Select()
